$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The G column (rows 2-27) previously held per-row sample codes (E7760..E7766);
# consolidate them all down to a single new code, E7420.
$ws.Range("G2:G27").Value = "E7420"
$ws.Range("G2:G27").Font.Name = "Arial"

# The H column (rows 2-27) held a literal boolean FALSE; turn each into a
# live =FALSE() formula instead.
$ws.Range("H2:H27").Formula = "=FALSE()"

# Selection cursor moved from the H column to the G column.
$ws.Range("G2:G27").Select()
